# Updated cryptos list: refreshed prices and 1h volume percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.353.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.09%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.347.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.71%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "410.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.10%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "114.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.52%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.589"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.60%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.642"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.49%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.16%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0994"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.86%  "

# Row 12
$ws.Range("E12").Value = "  +1.53%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.874.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.58%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.56%  "

# Row 15
$ws.Range("E15").Value = "  +1.46%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.358.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.99%  "

# Row 17
$ws.Range("E17").Value = "  +0.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "59.117.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.74%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.59%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000111"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.49%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "305.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.34%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.84%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.80%  "

# Row 27
$ws.Range("E27").Value = "  +2.78%  "

# Row 28
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.05%  "

# Row 29
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.179"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.58%  "

# Row 31
$ws.Range("E31").Value = "  +5.76%  "

# Row 32
$ws.Range("E32").Value = "  +0.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.87%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.86%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0521"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.89%  "

# Row 36
$ws.Range("E36").Value = "  +0.14%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.45%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.88%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.95%  "

# Row 41
$ws.Range("E41").Value = "  +3.12%  "

# Row 42
$ws.Range("E42").Value = "  +2.56%  "

# Row 43
$ws.Range("E43").Value = "  +0.12%  "

# Row 44
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.10%  "

# Row 45
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.52%  "

# Row 46
$ws.Range("E46").Value = "  -1.40%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.89%  "

# Row 48
$ws.Range("E48").Value = "  +7.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.219.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.26%  "

# Row 50
$ws.Range("E50").Value = "  -0.90%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.74%  "
